# Equipment formulario: switch the "Sucursales" lookup placeholders for
# "Equipos" ones, add a Num_Serie / Sucursal mini-table (rows 7-8 headers +
# a repeating body rows 9-41), rename the defined name, and drop the merges
# that no longer apply to the new single-column body.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Workbook-level defined name: Sucursales (A3:F19) -> Valores (A8:B9)
# ---------------------------------------------------------------------
$wb.Names.Item("Sucursales").Delete()
[void]$ws.Names.Add("Valores", "=Equipos!`$A`$8:`$B`$9")

# ---------------------------------------------------------------------
# 2) Re-point the existing placeholder cells from Sucursales.* to Equipos.*
# ---------------------------------------------------------------------
$ws.Range("B3").Value = "{{Equipos.Clave}}"
$ws.Range("E3").Value = "{{Equipos.CategoriaText}}"
$ws.Range("B5").Value = "{{Equipos.Nombre}}"
$ws.Range("E5").Value = "{{Equipos.Activo}}"

# ---------------------------------------------------------------------
# 3) New mini-table header (row 7) + item row (row 8): Num. Serie / Sucursal
# ---------------------------------------------------------------------
$ws.Range("A7").Value = "Número de Serie"
$ws.Range("B7").Value = "Sucursal"
$ws.Range("A8").Value = "{{item.Num_Serie}}"
$ws.Range("B8").Value = "{{item.SucursalText}}"

# ---------------------------------------------------------------------
# 4) Undo the merges that belonged to the old 2-column (B:C / E:F) layout
#    from row 7 down - the new body only needs single A / B columns.
# ---------------------------------------------------------------------
$oldMerges = @("B7:C7","E7:F7","B9:C9","E9:F9","B11:C11","E11:F11","B13:C13","B15:C15","B17:C17","B19:C19")
foreach ($rng in $oldMerges) {
    [void]$ws.Range($rng).UnMerge()
}

# ---------------------------------------------------------------------
# 5) Re-style row 7 (bold, centered mini-header) and row 8 (plain, centered)
# ---------------------------------------------------------------------
$ws.Range("A7:B7").Font.Bold = $true
$ws.Range("A7:B7").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A7:B7").VerticalAlignment = -4108     # xlCenter

$ws.Range("C7").HorizontalAlignment = -4142      # xlGeneral (clears explicit left)
$ws.Range("C7").VerticalAlignment = -4108        # xlCenter
$ws.Range("E7:F7").HorizontalAlignment = -4142
$ws.Range("E7:F7").VerticalAlignment = -4108

$ws.Range("A8").Font.Bold = $true
$ws.Range("A8").HorizontalAlignment = -4108
$ws.Range("A8").VerticalAlignment = -4108
$ws.Range("B8").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 6) Rows 9-41: repeating body of the Num_Serie / Sucursal mini-table.
#    Column A keeps the bold "header-like" look the template already used
#    for its left labels; column B is plain + centered; the former C/E/F
#    helper cells lose their left alignment (fall back to vertical-center
#    only, matching the old spacer style) wherever they still exist.
# ---------------------------------------------------------------------
for ($r = 9; $r -le 19; $r += 2) {
    $ws.Range("A$r").Font.Bold = $true
    $ws.Range("A$r").HorizontalAlignment = -4108
    $ws.Range("A$r").VerticalAlignment = -4108

    $ws.Range("B$r").HorizontalAlignment = -4108
    $ws.Range("B$r").VerticalAlignment = -4108

    if ($ws.Range("C$r").Value -ne $null -or $r -le 17) {
        $ws.Range("C$r").HorizontalAlignment = -4142
        $ws.Range("C$r").VerticalAlignment = -4108
    }
    if ($r -le 11) {
        $ws.Range("E$r").HorizontalAlignment = -4142
        $ws.Range("E$r").VerticalAlignment = -4108
        $ws.Range("F$r").HorizontalAlignment = -4142
        $ws.Range("F$r").VerticalAlignment = -4108
    }
}

for ($r = 10; $r -le 16; $r += 2) {
    $ws.Range("A$r").Font.Bold = $true
    $ws.Range("A$r").HorizontalAlignment = -4108
    $ws.Range("A$r").VerticalAlignment = -4108

    $ws.Range("B$r").HorizontalAlignment = -4108
}

# Rows 18-23: plain (non-bold) centered cells, both columns
for ($r = 18; $r -le 23; $r++) {
    $ws.Range("A$r").HorizontalAlignment = -4108
    $ws.Range("B$r").HorizontalAlignment = -4108
}

# Rows 24-41: brand-new blank body rows, plain centered cells
for ($r = 24; $r -le 41; $r++) {
    $ws.Range("A$r").Value = ""
    $ws.Range("B$r").Value = ""
    $ws.Range("A$r").HorizontalAlignment = -4108
    $ws.Range("B$r").HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------
# 7) Selection moves from D12 to A7
# ---------------------------------------------------------------------
[void]$ws.Range("A7").Select()
